$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# Reposition / resize the subtitle placeholder (EMU -> points, 12700 EMU per point)
$shp.Left = 611560 / 12700
$shp.Top = 3501008 / 12700
$shp.Width = 8280920 / 12700
$shp.Height = 1752600 / 12700

# Rename "Chương 6. Mô hình ngôn ngữ" -> "Bài 6. Mô hình ngôn ngữ", split across two runs
$tr = $shp.TextFrame.TextRange
$titlePara = $tr.Paragraphs(1, 1)
$titlePara.Text = "Bài 6. Mô hình ngôn ngữ"

$firstRun = $titlePara.Characters(1, 5)
$firstRun.Text = "Bài 6"
